{"js": "// Sequentially replace the date line and each \"A\u00f7B=\" math-fact cell with\n// its new value, in document order. A positional mapping is required\n// (rather than a global find/replace) because some \"before\" values\n// (e.g. \"90\u00f77=\") occur more than once in the document but map to\n// different \"after\" values depending on position.\nconst replacements = [\n  \"2026-01-01 Thursday\",\n  \"90\u00f73=\",\n  \"75\u00f72=\",\n  \"15\u00f79=\",\n  \"76\u00f74=\",\n  \"10\u00f73=\",\n  \"74\u00f78=\",\n  \"86\u00f73=\",\n  \"35\u00f75=\",\n  \"10\u00f76=\",\n  \"23\u00f75=\",\n  \"27\u00f72=\",\n  \"86\u00f78=\",\n  \"45\u00f76=\",\n  \"66\u00f74=\",\n  \"57\u00f75=\",\n  \"17\u00f79=\",\n  \"11\u00f74=\",\n  \"90\u00f72=\",\n  \"73\u00f73=\",\n  \"48\u00f73=\",\n  \"81\u00f79=\",\n  \"55\u00f77=\",\n  \"90\u00f75=\",\n  \"84\u00f75=\",\n  \"80\u00f78=\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet idx = 0;\nfor (let i = 0; i < paragraphs.items.length && idx < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text !== \"\") {\n    para.insertText(replacements[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and every \"A\u00f7B=\" math-fact cell in the table to\n# their new values. Cells are addressed by explicit (row, column)\n# coordinates rather than a blanket Find/Replace because some \"before\"\n# values (e.g. \"90\u00f77=\") appear more than once in the table but map to\n# different \"after\" values depending on which cell they are in.\n\n$d = $word.ActiveDocument\n\n# Date line (first paragraph, above the table).\n$d.Paragraphs(1).Range.Text = \"2026-01-01 Thursday\"\n\n$tbl = $d.Tables(1)\n\n# Row => ordered list of new values for columns 1..5.\n$newValues = @{\n    1  = @(\"90\u00f73=\", \"75\u00f72=\", \"15\u00f79=\", \"76\u00f74=\", \"10\u00f73=\")\n    5  = @(\"74\u00f78=\", \"86\u00f73=\", \"35\u00f75=\", \"10\u00f76=\", \"23\u00f75=\")\n    9  = @(\"27\u00f72=\", \"86\u00f78=\", \"45\u00f76=\", \"66\u00f74=\", \"57\u00f75=\")\n    13 = @(\"17\u00f79=\", \"11\u00f74=\", \"90\u00f72=\", \"73\u00f73=\", \"48\u00f73=\")\n    17 = @(\"81\u00f79=\", \"55\u00f77=\", \"90\u00f75=\", \"84\u00f75=\", \"80\u00f78=\")\n}\n\nforeach ($row in $newValues.Keys) {\n    $values = $newValues[$row]\n    for ($col = 1; $col -le 5; $col++) {\n        $tbl.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
